$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("A1").Value = "PRODUTO"
$ws.Range("B1").Value = "PREÇO"
$ws.Range("C1").Value = "DATA ATUAL"

# --- Remove the now-redundant D:E columns (old duplicate PREÇO / DATA ATUAL) ---
$ws.Range("D1:E4").Delete()

# --- Row 2: Dolar -> Dolar/BRL ---
$ws.Range("A2").Value = "Dolar/BRL"
$ws.Range("B2").Value = 4.81633799834
$ws.Range("C2").Value = 44704.81221411736

# --- Row 3: Café -> Café/Kg ---
$ws.Range("A3").Value = "Café/Kg"
$ws.Range("B3").Value = 2.161274399999995

# --- Row 4: Petroleo price updated ---
$ws.Range("B4").Value = 113.16292799859

# --- New rows 5-14 ---
$ws.Range("A5").Value = "Trigo"
$ws.Range("B5").Value = 422.3325130832395
$ws.Range("C5").Value = ""

$ws.Range("A6").Value = "Algodão"
$ws.Range("B6").Value = 1.430484650349657
$ws.Range("C6").Value = ""

$ws.Range("A7").Value = "Açucar"
$ws.Range("B7").Value = 0.2023615560060839
$ws.Range("C7").Value = ""

$ws.Range("A8").Value = "Arroz"
$ws.Range("B8").Value = 17.34174548872191
$ws.Range("C8").Value = ""

$ws.Range("A9").Value = "Etanol"
$ws.Range("B9").Value = 2.161999740740733
$ws.Range("C9").Value = ""

$ws.Range("A10").Value = "Feijao"
$ws.Range("B10").Value = 17.07222611736823
$ws.Range("C10").Value = ""

$ws.Range("A11").Value = "Gas Natural"
$ws.Range("B11").Value = 7.95004838783437
$ws.Range("C11").Value = ""

$ws.Range("A12").Value = "Madeira"
$ws.Range("B12").Value = 0.6758947314687541
$ws.Range("C12").Value = ""

$ws.Range("A13").Value = "Borracha"
$ws.Range("B13").Value = 142.8862158947376
$ws.Range("C13").Value = ""

$ws.Range("A14").Value = "Milho"
$ws.Range("B14").Value = 7.857316298887055
$ws.Range("C14").Value = ""

Write-Output "edit applied"
